$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new transaction was logged on 2024-09-07, pushing all existing
# September (and subsequently August) rows down by one row.
$ws.Rows.Item(35).Insert()

# Populate the freshly inserted row with the new September entry.
$ws.Range("R35").Value = "balance your axis"
$ws.Range("S35").Value = "2024-09-07 09:34:58"
